$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove unused hyperlink-related cell styles (left over from template; no cells use them)
foreach ($styleName in @("Hipervínculo", "Hipervínculo visitado")) {
    try {
        $wb.Styles.Item($styleName).Delete()
    } catch {
    }
}

# Replace the sample data (rows 2-18) with the new CRM pipeline data
# Row 2
$ws.Range("A2").Value = 714466
$ws.Range("B2").Value = 'Week Password'
$ws.Range("C2").Value = 'Craig Booker'
$ws.Range("D2").Value = 'Juan Para'
$ws.Range("E2").Value = 'CPU'
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 30000
$ws.Range("H2").Value = 'presented'

# Row 3
$ws.Range("A3").Value = 714466
$ws.Range("B3").Value = 'CBTS overdue'
$ws.Range("C3").Value = 'Craig Booker'
$ws.Range("D3").Value = 'Juan Para'
$ws.Range("E3").Value = 'Software'
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10000
$ws.Range("H3").Value = 'presented'

# Row 4
$ws.Range("A4").Value = 714466
$ws.Range("B4").Value = 'PhishMe clicked'
$ws.Range("C4").Value = 'Craig Booker'
$ws.Range("D4").Value = 'Juan Para'
$ws.Range("E4").Value = 'Maintenance'
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 5000
$ws.Range("H4").Value = 'pending'

# Row 5
$ws.Range("A5").Value = 737550
$ws.Range("B5").Value = 'Security Incident Involvement'
$ws.Range("C5").Value = 'Craig Booker'
$ws.Range("D5").Value = 'Juan Para'
$ws.Range("E5").Value = 'CPU'
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 35000
$ws.Range("H5").Value = 'declined'

# Row 6
$ws.Range("A6").Value = 146832
$ws.Range("B6").Value = 'Week Password'
$ws.Range("C6").Value = 'Daniel Hilton'
$ws.Range("D6").Value = 'Felipe Fiorin'
$ws.Range("E6").Value = 'CPU'
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 65000
$ws.Range("H6").Value = 'won'

# Row 7
$ws.Range("A7").Value = 218895
$ws.Range("B7").Value = 'CBTS overdue'
$ws.Range("C7").Value = 'Daniel Hilton'
$ws.Range("D7").Value = 'Felipe Fiorin'
$ws.Range("E7").Value = 'CPU'
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 40000
$ws.Range("H7").Value = 'pending'

# Row 8
$ws.Range("A8").Value = 218895
$ws.Range("B8").Value = 'PhishMe clicked'
$ws.Range("C8").Value = 'Daniel Hilton'
$ws.Range("D8").Value = 'Felipe Fiorin'
$ws.Range("E8").Value = 'Software'
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10000
$ws.Range("H8").Value = 'presented'

# Row 9
$ws.Range("A9").Value = 412290
$ws.Range("B9").Value = 'Security Incident Involvement'
$ws.Range("C9").Value = 'John Smith'
$ws.Range("D9").Value = 'Felipe Fiorin'
$ws.Range("E9").Value = 'Maintenance'
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 5000
$ws.Range("H9").Value = 'pending'

# Row 10
$ws.Range("A10").Value = 740150
$ws.Range("B10").Value = 'Week Password'
$ws.Range("C10").Value = 'John Smith'
$ws.Range("D10").Value = 'Kevin Whelan'
$ws.Range("E10").Value = 'CPU'
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 35000
$ws.Range("H10").Value = 'declined'

# Row 11
$ws.Range("A11").Value = 141962
$ws.Range("B11").Value = 'CBTS overdue'
$ws.Range("C11").Value = 'Cedric Moss'
$ws.Range("D11").Value = 'Kevin Whelan'
$ws.Range("E11").Value = 'CPU'
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 65000
$ws.Range("H11").Value = 'won'

# Row 12
$ws.Range("A12").Value = 163416
$ws.Range("B12").Value = 'PhishMe clicked'
$ws.Range("C12").Value = 'Cedric Moss'
$ws.Range("D12").Value = 'Kevin Whelan'
$ws.Range("E12").Value = 'CPU'
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 30000
$ws.Range("H12").Value = 'presented'

# Row 13
$ws.Range("A13").Value = 239344
$ws.Range("B13").Value = 'Security Incident Involvement'
$ws.Range("C13").Value = 'Cedric Moss'
$ws.Range("D13").Value = 'Kevin Whelan'
$ws.Range("E13").Value = 'Maintenance'
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 5000
$ws.Range("H13").Value = 'pending'

# Row 14
$ws.Range("A14").Value = 239344
$ws.Range("B14").Value = 'Week Password'
$ws.Range("C14").Value = 'Cedric Moss'
$ws.Range("D14").Value = 'Maeve Morris'
$ws.Range("E14").Value = 'Software'
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 10000
$ws.Range("H14").Value = 'presented'

# Row 15
$ws.Range("A15").Value = 307599
$ws.Range("B15").Value = 'CBTS overdue'
$ws.Range("C15").Value = 'Wendy Yule'
$ws.Range("D15").Value = 'Maeve Morris'
$ws.Range("E15").Value = 'Maintenance'
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 7000
$ws.Range("H15").Value = 'won'

# Row 16
$ws.Range("A16").Value = 688981
$ws.Range("B16").Value = 'PhishMe clicked'
$ws.Range("C16").Value = 'Wendy Yule'
$ws.Range("D16").Value = 'Maeve Morris'
$ws.Range("E16").Value = 'CPU'
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 100000
$ws.Range("H16").Value = 'won'

# Row 17
$ws.Range("A17").Value = 729833
$ws.Range("B17").Value = 'Security Incident Involvement'
$ws.Range("C17").Value = 'Wendy Yule'
$ws.Range("D17").Value = 'Maeve Morris'
$ws.Range("E17").Value = 'CPU'
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 65000
$ws.Range("H17").Value = 'declined'

# Row 18
$ws.Range("A18").Value = 729833
$ws.Range("B18").Value = 'Koepp Ltd'
$ws.Range("C18").Value = 'Wendy Yule'
$ws.Range("D18").Value = 'Fred Anderson'
$ws.Range("E18").Value = 'Monitor'
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 5000
$ws.Range("H18").Value = 'presented'

# Update the active selection to B1 (also clears the stale top-left scroll anchor)
$ws.Range("B1").Select()
